$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing style (style index 1, already applied to A2) so every
# newly-written cell keeps the same formatting as its row/column neighbours.
$ws.Range("A2").Copy()

$fillCells = @("B2","C2","D2","E2","F2","B3","C3","D3","E3","F3","C4","D4","E4","F4","D5","E5","F5","E6","F6","F7")
foreach ($addr in $fillCells) {
    $ws.Range($addr).PasteSpecial(-4122)
}

# Set the actual "-" text values for the newly formatted cells.
$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"
$ws.Range("F2").Value = "-"

$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = "-"
$ws.Range("F3").Value = "-"

$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "-"

$ws.Range("D5").Value = "-"
$ws.Range("E5").Value = "-"
$ws.Range("F5").Value = "-"

$ws.Range("E6").Value = "-"
$ws.Range("F6").Value = "-"

$ws.Range("F7").Value = "-"
